# Automatische test-sync: 2025-08-06 19:40:50
# Appends a new log row to the "Logs" sheet and updates the matching
# summary count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 6

$logs.Cells.Item($newRow, 1).Value = "Los jij dit even af?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Los jij dit even af?"
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-06 19:40:29"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Bump the Dashboard "Planning / Afspraak" count to include the new row.
$dashboard.Cells.Item(2, 2).Value = 5

# Extend the conditional-formatting ranges so the new row is covered too.
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "5")
    $newRange = $logs.Range($col + "2:" + $col + "6")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}
